# Revision of first 3 examples (DD-based recipes)
# Applies the example_01 input.xlsx changes:
#  - data fix: C3 2 -> ... (dd value correction)
#  - refreshed row heights (re-measured on re-save)
#  - moved active selection to F7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data correction -------------------------------------------------
# C3 ("dd" for the second recipe row) corrected from 1 to 2.
$ws.Range("C3").Value = 2

# --- Row height refresh (matches re-measured heights in the revision) -
$ws.Rows.Item(1).RowHeight = 51.75
$ws.Rows.Item(2).RowHeight = 39
$ws.Rows.Item(3).RowHeight = 51.75

# --- Selection moved to F7 --------------------------------------------
$ws.Range("F7").Select()
